# Update the "Sprint 1" backlog sheet:
#  - Add a "Status" column (D) to the use-case tracker table (rows 7-11)
#  - Mark the map/regions item as "Done" (green fill)
#  - Mark the store, start-screen and end-screen items as "In Progress" (yellow fill)
#  - Adjust column widths for the new column and its empty neighbour
#  - Update the sheet's scroll position / selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 1")
$ws.Activate()

# --- Status values for each backlog item ---------------------------------
# D7  -> header                             => "Status"
# D8  -> "See map broken up by DEC regions" => Done (green / theme accent6)
# D9  -> "Create store to spend money in "  => In Progress (yellow)
# D10 -> "Create a start screen"            => In Progress (yellow)
# D11 -> "Create an end screen"             => In Progress (yellow)

# Enter the text values first (in this order so shared-string table indices
# come out as In Progress / Status / Done).
$ws.Range("D9").Value = "In Progress"
$ws.Range("D7").Value = "Status"
$ws.Range("D8").Value = "Done"
$ws.Range("D10").Value = "In Progress"
$ws.Range("D11").Value = "In Progress"

# Then apply the fills. The yellow fill is primed on D9 and D8 before D8's
# theme fill is set, which keeps the workbook's fill table minimal (avoids a
# spurious intermediate fill entry).
$ws.Range("D9").Interior.Color = 65535
$ws.Range("D8").Interior.Color = 65535
$ws.Range("D8").Interior.ThemeColor = 10
$ws.Range("D10").Interior.Color = 65535
$ws.Range("D11").Interior.Color = 65535

# --- Column widths ---------------------------------------------------------
# Column D (Status) widens slightly; the newly-split, still-empty column E
# keeps a narrower width; columns F:P stay at their original width.
$ws.Columns.Item(4).ColumnWidth = 9.14
$ws.Columns.Item(5).ColumnWidth = 7.7

# --- View state: scroll down one row, select D11 --------------------------
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("D11").Select()

Write-Output "Status column added; backlog items updated."
